# chore: update Sheets via scheduled runner
# Refresh computed pricing/profit columns (H,I,J,K,L,M,N) across the
# per-craft "Leve" tables (ALC, ARM, BSM, CRP, CUL, LTW, WVR) with the
# latest market-board snapshot values.

$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 358.8
$ws.Range("I2").Value = 287.55554
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 287.55554
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -174.55554
$ws.Range("N2").Value = -1226
# Row 17
$ws.Range("H17").Value = 3097.9363
$ws.Range("J17").Value = 3097.9363
$ws.Range("L17").Value = 9293.8089
$ws.Range("N17").Value = -9629.8089
# Row 48
$ws.Range("H48").Value = 3999.8
$ws.Range("J48").Value = 3999.8
$ws.Range("L48").Value = 11999.4
$ws.Range("N48").Value = -12583.4
# Row 56
$ws.Range("H56").Value = 3999.8
$ws.Range("J56").Value = 3999.8
$ws.Range("L56").Value = 11999.4
$ws.Range("N56").Value = -13067.4
# Row 74
$ws.Range("H74").Value = 4071.4285
$ws.Range("J74").Value = 4187.5
$ws.Range("L74").Value = 4187.5
$ws.Range("N74").Value = -6059.5
# Row 77
$ws.Range("H77").Value = 4071.4285
$ws.Range("J77").Value = 4187.5
$ws.Range("L77").Value = 20937.5
$ws.Range("N77").Value = -30297.5
# Row 107
$ws.Range("H107").Value = 607.25
$ws.Range("I107").Value = 529.28
$ws.Range("J107").Value = 885.7143
$ws.Range("K107").Value = 529.28
$ws.Range("L107").Value = 885.7143
$ws.Range("M107").Value = 1390.72
$ws.Range("N107").Value = -4725.7143
# Row 112
$ws.Range("H112").Value = 1539.5333
$ws.Range("J112").Value = 1596.0358
$ws.Range("L112").Value = 4788.107400000001
$ws.Range("N112").Value = -7004.107400000001
# Row 132
$ws.Range("H132").Value = 2168.7446
$ws.Range("I132").Value = 1648.2307
$ws.Range("J132").Value = 4706.25
$ws.Range("K132").Value = 4944.6921
$ws.Range("L132").Value = 14118.75
$ws.Range("M132").Value = -2414.6921
$ws.Range("N132").Value = -19178.75
# Row 137
$ws.Range("H137").Value = 3648
$ws.Range("I137").Value = 2112.8
$ws.Range("J137").Value = 4694.727
$ws.Range("K137").Value = 6338.400000000001
$ws.Range("L137").Value = 14084.181
$ws.Range("M137").Value = -3788.400000000001
$ws.Range("N137").Value = -19184.181
# Row 138
$ws.Range("H138").Value = 1063318.5
$ws.Range("I138").Value = 2573.8667
$ws.Range("J138").Value = 1287419.4
$ws.Range("K138").Value = 7721.6001
$ws.Range("L138").Value = 3862258.2
$ws.Range("M138").Value = -2581.6001
$ws.Range("N138").Value = -3872538.2

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2475.3684
$ws.Range("I2").Value = 2501.7778
$ws.Range("K2").Value = 2501.7778
$ws.Range("M2").Value = -2388.7778
# Row 32
$ws.Range("H32").Value = 22846.623
$ws.Range("I32").Value = 11609.906
$ws.Range("J32").Value = 39969.24
$ws.Range("K32").Value = 11609.906
$ws.Range("L32").Value = 39969.24
$ws.Range("M32").Value = -11322.906
$ws.Range("N32").Value = -40543.24
# Row 97
$ws.Range("H97").Value = 1152.2069
$ws.Range("I97").Value = 834.913
$ws.Range("J97").Value = 2368.5
$ws.Range("K97").Value = 834.913
$ws.Range("L97").Value = 2368.5
$ws.Range("M97").Value = -338.913
$ws.Range("N97").Value = -3360.5
# Row 116
$ws.Range("H116").Value = 2475.3684
$ws.Range("I116").Value = 2501.7778
$ws.Range("K116").Value = 2501.7778
$ws.Range("M116").Value = -207.7777999999998
# Row 129
$ws.Range("H129").Value = 32848.5
$ws.Range("J129").Value = 33419.555
$ws.Range("L129").Value = 33419.555
$ws.Range("N129").Value = -43419.555

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2475.3684
$ws.Range("I3").Value = 2501.7778
$ws.Range("K3").Value = 2501.7778
$ws.Range("M3").Value = -2387.7778
# Row 105
$ws.Range("H105").Value = 7910.091
$ws.Range("I105").Value = 9200
$ws.Range("K105").Value = 9200
$ws.Range("M105").Value = -7453

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4633.727
$ws.Range("I31").Value = 2227.8333
$ws.Range("K31").Value = 2227.8333
$ws.Range("M31").Value = -1932.8333
# Row 34
$ws.Range("H34").Value = 4633.727
$ws.Range("I34").Value = 2227.8333
$ws.Range("K34").Value = 2227.8333
$ws.Range("M34").Value = -2025.8333
# Row 59
$ws.Range("H59").Value = 12000
$ws.Range("I59").Value = 12000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 12000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -10855
$ws.Range("N59").ClearContents()
# Row 76
$ws.Range("H76").Value = 9102.857
$ws.Range("I76").Value = 9102.857
$ws.Range("K76").Value = 9102.857
$ws.Range("M76").Value = -8787.857
# Row 79
$ws.Range("H79").Value = 9102.857
$ws.Range("I79").Value = 9102.857
$ws.Range("K79").Value = 9102.857
$ws.Range("M79").Value = -8010.857
# Row 99
$ws.Range("H99").Value = 2011.1111
$ws.Range("I99").Value = 2011.1111
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2011.1111
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -513.1111000000001
$ws.Range("N99").ClearContents()
# Row 126
$ws.Range("H126").Value = 2011.1111
$ws.Range("I126").Value = 2011.1111
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6033.3333
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3563.3333

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2528050
$ws.Range("I5").Value = 657.0303
$ws.Range("J5").Value = 5055443
$ws.Range("K5").Value = 1971.0909
$ws.Range("L5").Value = 15166329
$ws.Range("M5").Value = -1859.0909
$ws.Range("N5").Value = -15166553
# Row 129
$ws.Range("H129").Value = 3362.2222
$ws.Range("I129").Value = 4852
$ws.Range("J129").Value = 1500
$ws.Range("K129").Value = 14556
$ws.Range("L129").Value = 4500
$ws.Range("M129").Value = -9556
$ws.Range("N129").Value = -14500
# Row 135
$ws.Range("H135").Value = 2528050
$ws.Range("I135").Value = 657.0303
$ws.Range("J135").Value = 5055443
$ws.Range("K135").Value = 5913.2727
$ws.Range("L135").Value = 45498987
$ws.Range("M135").Value = -3378.2727
$ws.Range("N135").Value = -45504057
# Row 138
$ws.Range("H138").Value = 6151.9644
$ws.Range("I138").Value = 8937.691999999999
$ws.Range("J138").Value = 3737.6667
$ws.Range("K138").Value = 26813.076
$ws.Range("L138").Value = 11213.0001
$ws.Range("M138").Value = -21673.076
$ws.Range("N138").Value = -21493.0001
# Row 139
$ws.Range("H139").Value = 1567889.9
$ws.Range("I139").Value = 2937099
$ws.Range("J139").Value = 3079.3809
$ws.Range("K139").Value = 8811297
$ws.Range("L139").Value = 9238.1427
$ws.Range("M139").Value = -8806157
$ws.Range("N139").Value = -19518.1427
# Row 140
$ws.Range("H140").Value = 2092.2092
$ws.Range("I140").Value = 1576.1471
$ws.Range("K140").Value = 4728.4413
$ws.Range("M140").Value = 451.5587000000005
# Row 141
$ws.Range("H141").Value = 4303.609
$ws.Range("I141").Value = 3365.5557
$ws.Range("J141").Value = 4906.643
$ws.Range("K141").Value = 10096.6671
$ws.Range("L141").Value = 14719.929
$ws.Range("M141").Value = -4916.667099999999
$ws.Range("N141").Value = -25079.929

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 543.25
$ws.Range("J46").Value = 685
$ws.Range("L46").Value = 685
$ws.Range("N46").Value = -1061
# Row 60
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
# Row 114
$ws.Range("H114").Value = 5000
$ws.Range("J114").Value = 5000
$ws.Range("L114").Value = 5000
$ws.Range("N114").Value = -13678

